# Hatların sağ ve solda olma durumunun tespiti sağlandı
# Flip the sign of the column G measurements for rows 2-9 so that lines
# on the opposite side are represented with negative values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = -51.94
$ws.Range("G3").Value2 = -52.04
$ws.Range("G4").Value2 = -52.14
$ws.Range("G5").Value2 = -50.14
$ws.Range("G6").Value2 = -50.06
$ws.Range("G7").Value2 = -49.97
$ws.Range("G8").Value2 = -49.89
$ws.Range("G9").Value2 = -208.89
